$d = $word.ActiveDocument

$newLines = @(
    '<i class="fa-solid fa-bars"></i>',
    '<i class="fa-brands fa-facebook-f"></i>',
    '<i class="fa-brands fa-twitter"></i>',
    '<i class="fa-brands fa-instagram"></i>',
    '<i class="fa-brands fa-whatsapp"></i>'
)

$end = $d.Content.End
$rng = $d.Range($end, $end)

$insertText = "`r" + ($newLines -join "`r")
$rng.InsertAfter($insertText)
